$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.553.95'
$ws.Range('E2').Value = '  -3.42%  '
$ws.Range('D3').Value = '3.345.81'
$ws.Range('E3').Value = '  -2.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.46'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.89'
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.92'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.417'
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('D12').Value = '3.918.06'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '3.341.91'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('D17').Value = '60.592.53'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.26'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.54'
$ws.Range('E19').Value = '  -0.94%  '
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '376.15'
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.561'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.74'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = '3.493.06'
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('E26').Value = '  -5.59%  '
$ws.Range('E27').Value = '  -4.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  -4.28%  '
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.73'
$ws.Range('E32').Value = '  -3.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.86'
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.29'
$ws.Range('E34').Value = '  -3.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.56'
$ws.Range('E36').Value = '  -5.02%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.84'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '167.69'
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.99'
$ws.Range('E39').Value = '  -13.14%  '
$ws.Range('D40').Value = '3.382.48'
$ws.Range('E40').Value = '  -2.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0751'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.29'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.14'
$ws.Range('E44').Value = '  -3.58%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.61'
$ws.Range('E45').Value = '  -5.27%  '
$ws.Range('D46').Value = '2.460.02'
$ws.Range('E46').Value = '  -4.41%  '
$ws.Range('E47').Value = '  -3.77%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.36'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.818'
$ws.Range('E51').Value = '  +0.17%  '
